$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("missing_values")
$ws.Range("B7").Value = 40
$ws.Range("C7").Value = 0.30849915162733299
$ws.Range("B8").Value = 2378
$ws.Range("C8").Value = 18.34027456424495
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = 0.095877277085330767
$ws.Range("B19").Value = 171
$ws.Range("C19").Value = 2.7325023969319271

$ws = $wb.Worksheets.Item("profile_missing_values")
$ws.Range("B3").Value = 33.052985702270817
$ws.Range("C3").Value = 54.970760233918128
$ws.Range("B4").Value = 66.947014297729183
$ws.Range("C4").Value = 45.029239766081872
$ws.Range("B6").Value = 0.33641715727502103
$ws.Range("C6").Value = 0.58479532163742687
$ws.Range("B7").Value = 12.363330529857022
$ws.Range("B8").Value = 22.329688814129518
$ws.Range("B9").Value = 23.25483599663583
$ws.Range("C9").Value = 1.1695906432748537
$ws.Range("B10").Value = 20.647603027754418
$ws.Range("C10").Value = 5.2631578947368416
$ws.Range("B11").Value = 14.507989907485284
$ws.Range("C11").Value = 39.1812865497076
$ws.Range("B12").Value = 6.5601345668629101
$ws.Range("C12").Value = 53.801169590643269
$ws.Range("B14").Value = 24.516400336417156
$ws.Range("C14").Value = 38.011695906432749
$ws.Range("B15").Value = 32.506307821698904
$ws.Range("C15").Value = 15.789473684210526
$ws.Range("B16").Value = 29.100084104289319
$ws.Range("C16").Value = 9.9415204678362574
$ws.Range("B17").Value = 10.933557611438182
$ws.Range("C17").Value = 12.865497076023392
$ws.Range("B18").Value = 2.8174936921783011
$ws.Range("C18").Value = 22.807017543859647
$ws.Range("B19").Value = 0.12615643397813289
$ws.Range("C19").Value = 0.58479532163742687
$ws.Range("B21").Value = 0.29436501261564341
$ws.Range("B22").Value = 0.12615643397813289
$ws.Range("C22").Value = 0.58479532163742687
$ws.Range("B23").Value = 26.703111858704791
$ws.Range("C23").Value = 49.122807017543856
$ws.Range("B24").Value = 46.761984861227923
$ws.Range("C24").Value = 23.391812865497073
$ws.Range("B25").Value = 5.9714045416316228
$ws.Range("C25").Value = 8.1871345029239766
$ws.Range("B26").Value = 14.760302775441547
$ws.Range("C26").Value = 8.7719298245614024
$ws.Range("B27").Value = 0.67283431455004206
$ws.Range("B28").Value = 4.7098402018502945
$ws.Range("C28").Value = 9.9415204678362574
$ws.Range("B30").Value = 9.3355761143818334
$ws.Range("C30").Value = 8.7719298245614024
$ws.Range("B31").Value = 6.9386038687973093
$ws.Range("C31").Value = 2.9239766081871341
$ws.Range("B32").Value = 9.9663582842724967
$ws.Range("C32").Value = 2.9239766081871341
$ws.Range("B33").Value = 66.904962153069803
$ws.Range("C33").Value = 12.280701754385964
$ws.Range("B34").Value = 0.50462573591253157
$ws.Range("B35").Value = 3.7005887300252311
$ws.Range("B36").Value = 2.3128679562657695
$ws.Range("C36").Value = 1.7543859649122806
$ws.Range("B37").Value = 0.33641715727502103
$ws.Range("C37").Value = 71.345029239766077
$ws.Range("B39").Value = 18.040370058873005
$ws.Range("C39").Value = 5.8479532163742682
$ws.Range("B40").Value = 0.58873002523128681
$ws.Range("C40").Value = 1.1695906432748537
$ws.Range("B41").Value = 1.808242220353238
$ws.Range("C41").Value = 0.58479532163742687
$ws.Range("B42").Value = 1.3877207737594619
$ws.Range("C42").Value = 0.58479532163742687
$ws.Range("B43").Value = 5.6349873843566023
$ws.Range("C43").Value = 2.3391812865497075
$ws.Range("B44").Value = 20.142977291841884
$ws.Range("C44").Value = 1.7543859649122806
$ws.Range("B45").Value = 11.606391925988225
$ws.Range("C45").Value = 2.3391812865497075
$ws.Range("B46").Value = 3.280067283431455
$ws.Range("C46").Value = 0.58479532163742687
$ws.Range("B47").Value = 8.9571068124474351
$ws.Range("C47").Value = 7.0175438596491224
$ws.Range("B48").Value = 28.216989066442387
$ws.Range("C48").Value = 6.4327485380116958
$ws.Range("B49").Value = 0.33641715727502103
$ws.Range("C49").Value = 71.345029239766077
$ws.Range("B51").Value = 19.638351555929354
$ws.Range("C51").Value = 10.526315789473683
$ws.Range("B52").Value = 80.361648444070639
$ws.Range("C52").Value = 89.473684210526315

$ws = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")
$ws.Range("G3").Value = 2122992.6601368403
$ws.Range("L3").Value = 4268412

$ws = $wb.Worksheets.Item("labor_jubpenimp_stochastic_reg")
$ws.Range("G3").Value = 1541093.4671792437
